$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 131144266
$ws.Range("AC2").Value = "Ringhack, äldre, på en gran."
$ws.Range("AF2").Value = ""
$ws.Range("AJ2").Value = "gran"
$ws.Range("AK2").Value = "Picea abies"
$ws.Range("AM2").Value = "Trädstam på levande träd"
$ws.Range("AO2").Value = "Stem on living tree # Picea abies"
$ws.Range("B2").Value = 57884
$ws.Range("E2").Value = 100109
$ws.Range("F2").Value = "Tretåig hackspett"
$ws.Range("G2").Value = "Picoides tridactylus"
$ws.Range("H2").Value = "(Linnaeus, 1758)"
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = ""
$ws.Range("L2").Value = ""
$ws.Range("M2").Value = "äldre spår"
$ws.Range("Q2").Value = 503547
$ws.Range("R2").Value = 7036037

# Row 3
$ws.Range("A3").Value = 131143975
$ws.Range("AC3").Value = "På rönn."
$ws.Range("AF3").Value = ""
$ws.Range("AJ3").Value = "rönn"
$ws.Range("AK3").Value = "Sorbus aucuparia"
$ws.Range("AM3").Value = "Bark på levande träd"
$ws.Range("AO3").Value = "Bark on living woody plant # Sorbus aucuparia"
$ws.Range("B3").Value = 80379
$ws.Range("D3").Value = "LC"
$ws.Range("E3").Value = 6462
$ws.Range("F3").Value = "Stuplav"
$ws.Range("G3").Value = "Nephroma bellum"
$ws.Range("H3").Value = "(Spreng.) Tuck."
$ws.Range("J3").Value = ""
$ws.Range("K3").Value = "med apothecier"
$ws.Range("L3").Value = ""
$ws.Range("M3").Value = ""
$ws.Range("Q3").Value = 503441
$ws.Range("R3").Value = 7036004

# Row 4
$ws.Range("A4").Value = 131143973
$ws.Range("K4").Value = ""
$ws.Range("Q4").Value = 503418
$ws.Range("R4").Value = 7036017

# Row 5
$ws.Range("A5").Value = 131143977
$ws.Range("B5").Value = 80350
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 6458
$ws.Range("F5").Value = "Lunglav"
$ws.Range("G5").Value = "Lobaria pulmonaria"
$ws.Range("H5").Value = "(L.) Hoffm."
$ws.Range("K5").Value = "med soral"
$ws.Range("Q5").Value = 503254
$ws.Range("R5").Value = 7036092

# Row 22
$ws.Range("A22").Value = 131143972
$ws.Range("AC22").Value = "På rönn."
$ws.Range("AJ22").Value = "rönn"
$ws.Range("AK22").Value = "Sorbus aucuparia"
$ws.Range("AM22").Value = "Bark på levande träd"
$ws.Range("AO22").Value = "Bark on living woody plant # Sorbus aucuparia"
$ws.Range("B22").Value = 80379
$ws.Range("D22").Value = "LC"
$ws.Range("E22").Value = 6462
$ws.Range("F22").Value = "Stuplav"
$ws.Range("G22").Value = "Nephroma bellum"
$ws.Range("H22").Value = "(Spreng.) Tuck."
$ws.Range("K22").Value = "med apothecier"
$ws.Range("Q22").Value = 503377
$ws.Range("R22").Value = 7036016

# Row 23
$ws.Range("A23").Value = 131143976
$ws.Range("AC23").Value = ""
$ws.Range("B23").Value = 80350
$ws.Range("D23").Value = "NT"
$ws.Range("E23").Value = 6458
$ws.Range("F23").Value = "Lunglav"
$ws.Range("G23").Value = "Lobaria pulmonaria"
$ws.Range("H23").Value = "(L.) Hoffm."
$ws.Range("K23").Value = ""
$ws.Range("Q23").Value = 503419
$ws.Range("R23").Value = 7036154

# Row 24
$ws.Range("A24").Value = 131143989
$ws.Range("AC24").Value = "Lunglav på en gran vid en rönn med lunglav."
$ws.Range("AJ24").Value = "gran"
$ws.Range("AK24").Value = "Picea abies"
$ws.Range("AM24").Value = "Gren på levande träd"
$ws.Range("AO24").Value = "Branch on living tree # Picea abies"
$ws.Range("Q24").Value = 503448
$ws.Range("R24").Value = 7036030

# Row 25
$ws.Range("A25").Value = 131143998
$ws.Range("AC25").Value = "På flera granar."
$ws.Range("AM25").Value = ""
$ws.Range("AO25").Value = "Picea abies"
$ws.Range("B25").Value = 79245
$ws.Range("E25").Value = 6425
$ws.Range("F25").Value = "Garnlav"
$ws.Range("G25").Value = "Alectoria sarmentosa"
$ws.Range("H25").Value = "(Ach.) Ach."
$ws.Range("Q25").Value = 503444
$ws.Range("R25").Value = 7036006

# Row 30
$ws.Range("A30").Value = 131144281
$ws.Range("AC30").Value = "På rönn."
$ws.Range("AJ30").Value = "rönn"
$ws.Range("AK30").Value = "Sorbus aucuparia"
$ws.Range("AM30").Value = "Bark på levande träd"
$ws.Range("AO30").Value = "Bark on living woody plant # Sorbus aucuparia"
$ws.Range("B30").Value = 80350
$ws.Range("E30").Value = 6458
$ws.Range("F30").Value = "Lunglav"
$ws.Range("G30").Value = "Lobaria pulmonaria"
$ws.Range("H30").Value = "(L.) Hoffm."
$ws.Range("Q30").Value = 503194
$ws.Range("R30").Value = 7035963

# Row 31
$ws.Range("A31").Value = 131143993
$ws.Range("Q31").Value = 503307
$ws.Range("R31").Value = 7036124

# Row 32
$ws.Range("A32").Value = 131143994
$ws.Range("AC32").Value = "På flera granar."
$ws.Range("AJ32").Value = "gran"
$ws.Range("AK32").Value = "Picea abies"
$ws.Range("AM32").Value = ""
$ws.Range("AO32").Value = "Picea abies"
$ws.Range("B32").Value = 79245
$ws.Range("E32").Value = 6425
$ws.Range("F32").Value = "Garnlav"
$ws.Range("G32").Value = "Alectoria sarmentosa"
$ws.Range("H32").Value = "(Ach.) Ach."
$ws.Range("Q32").Value = 503252
$ws.Range("R32").Value = 7036082

# Row 33
$ws.Range("A33").Value = 131144293
$ws.Range("Q33").Value = 503484
$ws.Range("R33").Value = 7036019

# Row 75
$ws.Range("A75").Value = 131144302
$ws.Range("AJ75").Value = "gran"
$ws.Range("AK75").Value = "Picea abies"
$ws.Range("AM75").Value = ""
$ws.Range("AO75").Value = "Picea abies"
$ws.Range("B75").Value = 79245
$ws.Range("D75").Value = "NT"
$ws.Range("E75").Value = 6425
$ws.Range("F75").Value = "Garnlav"
$ws.Range("G75").Value = "Alectoria sarmentosa"
$ws.Range("H75").Value = "(Ach.) Ach."
$ws.Range("K75").Value = ""
$ws.Range("Q75").Value = 503606
$ws.Range("R75").Value = 7036064

# Row 76
$ws.Range("A76").Value = 131144273
$ws.Range("AC76").Value = "På rönn."
$ws.Range("AJ76").Value = "rönn"
$ws.Range("AK76").Value = "Sorbus aucuparia"
$ws.Range("AM76").Value = "Bark på levande träd"
$ws.Range("AO76").Value = "Bark on living woody plant # Sorbus aucuparia"
$ws.Range("B76").Value = 80350
$ws.Range("E76").Value = 6458
$ws.Range("F76").Value = "Lunglav"
$ws.Range("G76").Value = "Lobaria pulmonaria"
$ws.Range("H76").Value = "(L.) Hoffm."
$ws.Range("Q76").Value = 503267
$ws.Range("R76").Value = 7036145

# Row 77
$ws.Range("A77").Value = 131143970
$ws.Range("AC77").Value = ""
$ws.Range("B77").Value = 80379
$ws.Range("D77").Value = "LC"
$ws.Range("E77").Value = 6462
$ws.Range("F77").Value = "Stuplav"
$ws.Range("G77").Value = "Nephroma bellum"
$ws.Range("H77").Value = "(Spreng.) Tuck."
$ws.Range("K77").Value = "med apothecier"
$ws.Range("Q77").Value = 503376
$ws.Range("R77").Value = 7035991
